$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 232.5
$ws.Range("I11").Value = 232.5
$ws.Range("K11").Value = 232.5
$ws.Range("M11").Value = -92.5

$ws.Range("H19").Value = 2298
$ws.Range("I19").Value = 949.75
$ws.Range("K19").Value = 949.75
$ws.Range("M19").Value = -774.75

$ws.Range("H43").Value = 6398.75
$ws.Range("I43").Value = 4800
$ws.Range("J43").Value = 6931.6665
$ws.Range("K43").Value = 4800
$ws.Range("L43").Value = 6931.6665
$ws.Range("M43").Value = -4731
$ws.Range("N43").Value = -7069.6665

$ws.Range("H70").Value = 2476.6316
$ws.Range("I70").Value = 2131.3333
$ws.Range("J70").Value = 2541.375
$ws.Range("K70").Value = 6393.999899999999
$ws.Range("L70").Value = 7624.125
$ws.Range("M70").Value = -6123.999899999999
$ws.Range("N70").Value = -8164.125

$ws.Range("H73").Value = 2476.6316
$ws.Range("I73").Value = 2131.3333
$ws.Range("J73").Value = 2541.375
$ws.Range("K73").Value = 6393.999899999999
$ws.Range("L73").Value = 7624.125
$ws.Range("M73").Value = -5457.999899999999
$ws.Range("N73").Value = -9496.125

$ws.Range("H106").Value = 2970.4546
$ws.Range("I106").Value = 2975
$ws.Range("K106").Value = 2975
$ws.Range("M106").Value = -2344

$ws.Range("H132").Value = 2238.2307
$ws.Range("I132").Value = 1908.8636
$ws.Range("K132").Value = 5726.5908
$ws.Range("M132").Value = -3196.5908

$ws.Range("H137").Value = 3284.7
$ws.Range("I137").Value = 2474.5
$ws.Range("K137").Value = 7423.5
$ws.Range("M137").Value = -4873.5

$ws.Range("H138").Value = 16951344
$ws.Range("I138").Value = 1156.2858
$ws.Range("J138").Value = 41670370
$ws.Range("K138").Value = 3468.8574
$ws.Range("L138").Value = 125011110
$ws.Range("M138").Value = 1671.1426
$ws.Range("N138").Value = -125021390


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15126.288
$ws.Range("I32").Value = 4501.6445
$ws.Range("K32").Value = 4501.6445
$ws.Range("M32").Value = -4214.6445

$ws.Range("H102").Value = 2640.6667
$ws.Range("J102").Value = 3831
$ws.Range("L102").Value = 3831
$ws.Range("N102").Value = -7075


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 12639.167
$ws.Range("J81").Value = 12639.167
$ws.Range("L81").Value = 12639.167
$ws.Range("N81").Value = -14761.167

$ws.Range("H84").Value = 12639.167
$ws.Range("J84").Value = 12639.167
$ws.Range("L84").Value = 37917.501
$ws.Range("N84").Value = -48525.501

$ws.Range("H86").Value = 2979.2
$ws.Range("I86").Value = 3043.7778
$ws.Range("J86").Value = 2398
$ws.Range("K86").Value = 3043.7778
$ws.Range("L86").Value = 2398
$ws.Range("M86").Value = -1920.7778
$ws.Range("N86").Value = -4644

$ws.Range("H89").Value = 2979.2
$ws.Range("I89").Value = 3043.7778
$ws.Range("J89").Value = 2398
$ws.Range("K89").Value = 15218.889
$ws.Range("L89").Value = 11990
$ws.Range("M89").Value = -9602.888999999999
$ws.Range("N89").Value = -23222

$ws.Range("H94").Value = 2942.6875
$ws.Range("I94").Value = 2363.0715
$ws.Range("K94").Value = 2363.0715
$ws.Range("M94").Value = -1912.0715


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 54.666668
$ws.Range("I7").Value = 37.38889
$ws.Range("K7").Value = 37.38889
$ws.Range("M7").Value = 75.61111

$ws.Range("H22").Value = 806.3333
$ws.Range("J22").Value = 1157.4445
$ws.Range("L22").Value = 1157.4445
$ws.Range("N22").Value = -1857.4445

$ws.Range("H58").Value = 6025.952
$ws.Range("I58").Value = 3828.75
$ws.Range("J58").Value = 8955.556
$ws.Range("K58").Value = 3828.75
$ws.Range("L58").Value = 8955.556
$ws.Range("M58").Value = -3625.75
$ws.Range("N58").Value = -9361.556

$ws.Range("H62").Value = 8198.6
$ws.Range("I62").Value = 7747.5
$ws.Range("K62").Value = 7747.5
$ws.Range("M62").Value = -7123.5

$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496

$ws.Range("H65").Value = 8198.6
$ws.Range("I65").Value = 7747.5
$ws.Range("K65").Value = 38737.5
$ws.Range("M65").Value = -35617.5

$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716

$ws.Range("H132").Value = 3187.1428
$ws.Range("I132").Value = 3281.64
$ws.Range("J132").Value = 2399.6667
$ws.Range("K132").Value = 9844.92
$ws.Range("L132").Value = 7199.000100000001
$ws.Range("M132").Value = -7314.92
$ws.Range("N132").Value = -12259.0001

$ws.Range("H134").Value = 5155.0347
$ws.Range("I134").Value = 4123.7144
$ws.Range("K134").Value = 12371.1432
$ws.Range("M134").Value = -9836.143199999999

$ws.Range("H136").Value = 6025.952
$ws.Range("I136").Value = 3828.75
$ws.Range("J136").Value = 8955.556
$ws.Range("K136").Value = 11486.25
$ws.Range("L136").Value = 26866.668
$ws.Range("M136").Value = -8936.25
$ws.Range("N136").Value = -31966.668


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 607
$ws.Range("I3").Value = 566.5
$ws.Range("K3").Value = 1699.5
$ws.Range("M3").Value = -1587.5

$ws.Range("H134").Value = 926.0625
$ws.Range("I134").Value = 754.4666999999999
$ws.Range("K134").Value = 2263.4001
$ws.Range("M134").Value = 2806.5999

$ws.Range("H139").Value = 2697.45
$ws.Range("I139").Value = 2430.9333
$ws.Range("K139").Value = 7292.7999
$ws.Range("M139").Value = -2152.7999

$ws.Range("H140").Value = 522258.5
$ws.Range("I140").Value = 1119
$ws.Range("K140").Value = 3357
$ws.Range("M140").Value = 1823


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 36333.332
$ws.Range("J69").Value = 36333.332
$ws.Range("L69").Value = 36333.332
$ws.Range("N69").Value = -37831.332

$ws.Range("H72").Value = 36333.332
$ws.Range("J72").Value = 36333.332
$ws.Range("L72").Value = 108999.996
$ws.Range("N72").Value = -116487.996


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 71666.5
$ws.Range("J20").Value = 71666.5
$ws.Range("L20").Value = 71666.5
$ws.Range("N20").Value = -72118.5

$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("K29").Value = 1000
$ws.Range("M29").Value = -705

$ws.Range("H55").Value = 662.2143
$ws.Range("I55").Value = 670.1111
$ws.Range("J55").Value = 648
$ws.Range("K55").Value = 670.1111
$ws.Range("L55").Value = 648
$ws.Range("M55").Value = -497.1111
$ws.Range("N55").Value = -994

$ws.Range("H82").Value = 994.3333
$ws.Range("I82").Value = 994
$ws.Range("J82").Value = 994.5
$ws.Range("K82").Value = 994
$ws.Range("L82").Value = 994.5
$ws.Range("M82").Value = -633
$ws.Range("N82").Value = -1716.5

$ws.Range("H85").Value = 994.3333
$ws.Range("I85").Value = 994
$ws.Range("J85").Value = 994.5
$ws.Range("K85").Value = 994
$ws.Range("L85").Value = 994.5
$ws.Range("M85").Value = 254
$ws.Range("N85").Value = -3490.5

$ws.Range("H136").Value = 2944.6428
$ws.Range("I136").Value = 2954.4602
$ws.Range("K136").Value = 8863.3806
$ws.Range("M136").Value = -6313.3806


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7494.5
$ws.Range("J15").Value = 7499
$ws.Range("L15").Value = 7499
$ws.Range("N15").Value = -8075

$ws.Range("H34").Value = 33333.332
$ws.Range("J34").Value = 20000
$ws.Range("L34").Value = 20000
$ws.Range("N34").Value = -20406

$ws.Range("H41").Value = 19378.545
$ws.Range("I41").Value = 13000
$ws.Range("J41").Value = 20016.4
$ws.Range("K41").Value = 13000
$ws.Range("L41").Value = 20016.4
$ws.Range("M41").Value = -12610
$ws.Range("N41").Value = -20796.4

$ws.Range("H81").Value = 872
$ws.Range("I81").Value = 870.2353000000001
$ws.Range("K81").Value = 1740.4706
$ws.Range("M81").Value = -679.4706000000001

$ws.Range("H84").Value = 872
$ws.Range("I84").Value = 870.2353000000001
$ws.Range("K84").Value = 8702.353000000001
$ws.Range("M84").Value = -3398.353000000001

$ws.Range("H132").Value = 2699.0625
$ws.Range("I132").Value = 2824.5122
$ws.Range("J132").Value = 1964.2858
$ws.Range("K132").Value = 8473.536599999999
$ws.Range("L132").Value = 5892.857400000001
$ws.Range("M132").Value = -5943.536599999999
$ws.Range("N132").Value = -10952.8574

